# Auto-generated update of leve-profit calculation columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, matching refreshed
# Universalis market-board price data pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51: A Bile Business
$ws.Range("H51").Value = 3787.125
$ws.Range("I51").Value = 3499.3333
$ws.Range("K51").Value = 3499.3333
$ws.Range("M51").Value = -3015.3333

# Row 86: Filling in the Blanks
$ws.Range("H86").Value = 2528.8096
$ws.Range("I86").Value = 2075.2
$ws.Range("J86").Value = 3662.8333
$ws.Range("K86").Value = 2075.2
$ws.Range("L86").Value = 3662.8333
$ws.Range("M86").Value = -952.1999999999998
$ws.Range("N86").Value = -5908.8333

# Row 89: Ink into Antiquity (L)
$ws.Range("H89").Value = 2528.8096
$ws.Range("I89").Value = 2075.2
$ws.Range("J89").Value = 3662.8333
$ws.Range("K89").Value = 10376
$ws.Range("L89").Value = 18314.1665
$ws.Range("M89").Value = -4760
$ws.Range("N89").Value = -29546.1665

# Row 107: Another Man's Ink
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").ClearContents()
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = 0

# Row 113: Amaro Kart
$ws.Range("H113").Value = 4573.143
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

# Row 138: All-night Crafting
$ws.Range("H138").Value = 41668348
$ws.Range("I138").Value = 1426.7
$ws.Range("J138").Value = 71430440
$ws.Range("K138").Value = 4280.1
$ws.Range("L138").Value = 214291320
$ws.Range("M138").Value = 859.8999999999996
$ws.Range("N138").Value = -214301600

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 6688.524
$ws.Range("I32").Value = 3204.0588
$ws.Range("J32").Value = 21497.5
$ws.Range("K32").Value = 3204.0588
$ws.Range("L32").Value = 21497.5
$ws.Range("M32").Value = -2917.0588
$ws.Range("N32").Value = -22071.5

# Row 76: Sometimes the South Wins
$ws.Range("H76").Value = 141498
$ws.Range("J76").Value = 141498
$ws.Range("L76").Value = 141498
$ws.Range("N76").Value = -142174

# Row 79: The Thriller of Autumn (L)
$ws.Range("H79").Value = 141498
$ws.Range("J79").Value = 141498
$ws.Range("L79").Value = 141498
$ws.Range("N79").Value = -143838

# Row 97: Ore for Me
$ws.Range("H97").Value = 1097.7646
$ws.Range("I97").Value = 1036.9231
$ws.Range("K97").Value = 1036.9231
$ws.Range("M97").Value = -540.9231

# Row 103: Sweeping the Legs
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("N103").Value = 0

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 1717.375
$ws.Range("I132").Value = 1685.3556
$ws.Range("K132").Value = 5056.066800000001
$ws.Range("M132").Value = -2526.066800000001

$ws = $wb.Worksheets.Item("BSM")
# Row 80: Unbreaker
$ws.Range("H80").Value = 53020.79
$ws.Range("J80").Value = 462.15384
$ws.Range("L80").Value = 462.15384
$ws.Range("N80").Value = -2458.15384

# Row 83: Attack on Titanium (L)
$ws.Range("H83").Value = 53020.79
$ws.Range("J83").Value = 462.15384
$ws.Range("L83").Value = 2310.7692
$ws.Range("N83").Value = -12294.7692

# Row 94: High Steal
$ws.Range("H94").Value = 1064.875
$ws.Range("I94").Value = 1019.7778
$ws.Range("J94").Value = 1122.8572
$ws.Range("K94").Value = 1019.7778
$ws.Range("L94").Value = 1122.8572
$ws.Range("M94").Value = -568.7778
$ws.Range("N94").Value = -2024.8572

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 52953.75
$ws.Range("I105").Value = 168928.33
$ws.Range("J105").Value = 3250.3572
$ws.Range("K105").Value = 168928.33
$ws.Range("L105").Value = 3250.3572
$ws.Range("M105").Value = -167181.33
$ws.Range("N105").Value = -6744.3572

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 5101.5386
$ws.Range("I31").Value = 1825.75
$ws.Range("J31").Value = 10342.8
$ws.Range("K31").Value = 1825.75
$ws.Range("L31").Value = 10342.8
$ws.Range("M31").Value = -1530.75
$ws.Range("N31").Value = -10932.8

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 5101.5386
$ws.Range("I34").Value = 1825.75
$ws.Range("J34").Value = 10342.8
$ws.Range("K34").Value = 1825.75
$ws.Range("L34").Value = 10342.8
$ws.Range("M34").Value = -1623.75
$ws.Range("N34").Value = -10746.8

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 2532.9167
$ws.Range("I58").Value = 1743.375
$ws.Range("K58").Value = 1743.375
$ws.Range("M58").Value = -1540.375

# Row 62: Splinter in the Sewers
$ws.Range("H62").Value = 3525.25
$ws.Range("I62").Value = 3386.889
$ws.Range("J62").Value = 3703.1428
$ws.Range("K62").Value = 3386.889
$ws.Range("L62").Value = 3703.1428
$ws.Range("M62").Value = -2762.889
$ws.Range("N62").Value = -4951.1428

# Row 65: The Lumber of Their Discontent (L)
$ws.Range("H65").Value = 3525.25
$ws.Range("I65").Value = 3386.889
$ws.Range("J65").Value = 3703.1428
$ws.Range("K65").Value = 16934.445
$ws.Range("L65").Value = 18515.714
$ws.Range("M65").Value = -13814.445
$ws.Range("N65").Value = -24755.714

# Row 70: A Reward Fitting of the Faithful
$ws.Range("H70").Value = 55000
$ws.Range("J70").Value = 55000
$ws.Range("L70").Value = 55000
$ws.Range("N70").Value = -55630

# Row 73: Just Rewards for Just Devotion (L)
$ws.Range("H73").Value = 55000
$ws.Range("J73").Value = 55000
$ws.Range("L73").Value = 55000
$ws.Range("N73").Value = -57184

# Row 136: Turali Quality
$ws.Range("H136").Value = 2532.9167
$ws.Range("I136").Value = 1743.375
$ws.Range("K136").Value = 5230.125
$ws.Range("M136").Value = -2680.125

$ws = $wb.Worksheets.Item("CUL")
# Row 99: A Shorlonging for the Familiar
$ws.Range("H99").Value = 6012.5
$ws.Range("I99").Value = 25
$ws.Range("J99").Value = 12000
$ws.Range("K99").Value = 75
$ws.Range("L99").Value = 36000
$ws.Range("M99").Value = 2171
$ws.Range("N99").Value = -40492

# Row 122: Salt of the North
$ws.Range("H122").Value = 594758.9
$ws.Range("J122").Value = 1010708.8
$ws.Range("L122").Value = 9096379.200000001
$ws.Range("N122").Value = -9101279.200000001

$ws = $wb.Worksheets.Item("GSM")
# Row 52: It's My Business to Know Things
$ws.Range("H52").Value = 20000.2
$ws.Range("J52").Value = 20000.2
$ws.Range("L52").Value = 20000.2
$ws.Range("N52").Value = -20518.2

# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 41669150
$ws.Range("I80").Value = 76925440
$ws.Range("J80").Value = 2628.3635
$ws.Range("K80").Value = 76925440
$ws.Range("L80").Value = 2628.3635
$ws.Range("M80").Value = -76924442
$ws.Range("N80").Value = -4624.363499999999

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 41669150
$ws.Range("I83").Value = 76925440
$ws.Range("J83").Value = 2628.3635
$ws.Range("K83").Value = 384627200
$ws.Range("L83").Value = 13141.8175
$ws.Range("M83").Value = -384622208
$ws.Range("N83").Value = -23125.8175

# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 3031847
$ws.Range("I113").Value = 1240.25
$ws.Range("J113").Value = 4763622.5
$ws.Range("K113").Value = 1240.25
$ws.Range("L113").Value = 4763622.5
$ws.Range("M113").Value = 929.75
$ws.Range("N113").Value = -4767962.5

# Row 132: On Board for Lar
$ws.Range("H132").Value = 4468.514
$ws.Range("I132").Value = 3827.2727
$ws.Range("J132").Value = 5553.6924
$ws.Range("K132").Value = 11481.8181
$ws.Range("L132").Value = 16661.0772
$ws.Range("M132").Value = -8951.8181
$ws.Range("N132").Value = -21721.0772

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore
$ws.Range("H16").Value = 2357.4
$ws.Range("I16").Value = 2357.4
$ws.Range("K16").Value = 2357.4
$ws.Range("M16").Value = -2187.4

# Row 40: Best Served Toad
$ws.Range("H40").Value = 5558234.5
$ws.Range("I40").Value = 2976.5557
$ws.Range("K40").Value = 2976.5557
$ws.Range("M40").Value = -2840.5557

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 1682.5217
$ws.Range("I93").Value = 1603.5333
$ws.Range("K93").Value = 1603.5333
$ws.Range("M93").Value = -355.5333000000001

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 4487.326
$ws.Range("I136").Value = 3780.111
$ws.Range("J136").Value = 5492.316
$ws.Range("K136").Value = 11340.333
$ws.Range("L136").Value = 16476.948
$ws.Range("M136").Value = -8790.332999999999
$ws.Range("N136").Value = -21576.948

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire
$ws.Range("H122").Value = 2821.6943
$ws.Range("I122").Value = 1612.3334
$ws.Range("K122").Value = 4837.0002
$ws.Range("M122").Value = -2387.0002

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 18662.666
$ws.Range("I126").Value = 2158.6667
$ws.Range("K126").Value = 6476.000100000001
$ws.Range("M126").Value = -4006.000100000001
